# Appends a new observation record as row 35 of the "Artfynd" sheet
# (mirrors the shape of the existing rows 2-34).

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$row = 35

function Set-TextCell($cell, [string]$text) {
    # Force text storage so Excel's type-inference can't turn a
    # numeric-looking or date-looking string into a number/date,
    # then drop back to the default "Normal" style so no stray
    # number-format sticks to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-EmptyCell($cell) {
    # Materialize a present-but-empty cell (matches source rows that
    # carry blank values for these columns instead of omitting them).
    $cell.NumberFormat = "@"
    $cell.Style = "Normal"
}

# -- identifiers / numeric fields -----------------------------------
$ws.Cells.Item($row, 1).Value  = 111966065            # A  Id
$ws.Cells.Item($row, 2).Value  = 83148                # B  Taxonsorteringsordning
$ws.Cells.Item($row, 3).Value  = "Ovaliderad"          # C  Valideringsstatus
$ws.Cells.Item($row, 4).Value  = "LC"                  # D  Rödlistade
$ws.Cells.Item($row, 5).Value  = 3518                  # E  TaxonId
$ws.Cells.Item($row, 6).Value  = "Smal svampklubba"    # F  Artnamn
$ws.Cells.Item($row, 7).Value  = "Tolypocladium ophioglossoides"  # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "(Ehrh. ex J.F.Gmel.:Fr.) Quandt, Kepler & Spatafora"  # H  Auktor

Set-TextCell $ws.Cells.Item($row, 9) "2"               # I  Antal (stored as text)
$ws.Cells.Item($row, 10).Value = "fruktkroppar"        # J  Enhet

Set-EmptyCell $ws.Cells.Item($row, 11)                 # K  Ålder-Stadium (blank)
Set-EmptyCell $ws.Cells.Item($row, 14)                 # N  Metod (blank)

$ws.Cells.Item($row, 16).Value = "Angertuvan, öster om, Vg"  # P  Lokalnamn
$ws.Cells.Item($row, 17).Value = 338285.5070198396     # Q  Ost
$ws.Cells.Item($row, 18).Value = 6433442.904015562     # R  Nord
$ws.Cells.Item($row, 19).Value = 5                     # S  Noggrannhet
$ws.Cells.Item($row, 20).Value = "Västra Götaland"     # T  Län
$ws.Cells.Item($row, 21).Value = "Ale"                 # U  Kommun
$ws.Cells.Item($row, 22).Value = "Västergötland"       # V  Provins
$ws.Cells.Item($row, 23).Value = "Skepplanda"          # W  Församling

Set-TextCell $ws.Cells.Item($row, 25) "2023-09-06"     # Y  Startdatum
$ws.Cells.Item($row, 26).Value = "00:00"               # Z  Starttid
Set-TextCell $ws.Cells.Item($row, 27) "2023-09-06"     # AA Slutdatum
$ws.Cells.Item($row, 28).Value = "00:00"               # AB Sluttid

$ws.Cells.Item($row, 29).Value = "Växte vid största stigen som går vid foten av Angertuvans östra sluttning."  # AC Publik kommentar

$ws.Cells.Item($row, 30).Value = $false                # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false                # AE Osäker artbestämning
Set-EmptyCell $ws.Cells.Item($row, 32)                 # AF Bestämningsmetod (blank)
$ws.Cells.Item($row, 33).Value = $false                # AG Ospontan

Set-EmptyCell $ws.Cells.Item($row, 46)                 # AT Bestämningsår (blank)

$ws.Cells.Item($row, 49).Value = "Thomas Grönlund"     # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Thomas Grönlund"     # AX Observatörer

Set-EmptyCell $ws.Cells.Item($row, 51)                 # AY Projektnamn (blank)
